# Updated cryptos list on Sun Jul 23 10:02:30 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain text (e.g. "1.000", "0.9997",
# "29.866.45") even though many values look numeric. Force each such
# cell to a text number-format immediately before writing so Excel does
# not silently coerce it into a real number (which would drop trailing
# zeros, switch to scientific notation, or mis-parse "."-grouped values).
# Cells whose new value already contains two "." (e.g. "29.882.94")
# can never parse as a number, so they're written directly.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

$ws.Range("D2").Value = "29.882.94"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.873.72"
$ws.Range("E3").Value = "  -0.99%  "
Set-TextValue "D4" "0.9997"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  -4.40%  "
Set-TextValue "D6" "241.99"
$ws.Range("E6").Value = "  -0.62%  "
Set-TextValue "D7" "0.9991"
$ws.Range("E7").Value = "  -0.24%  "
Set-TextValue "D8" "0.3149"
$ws.Range("E8").Value = "  +0.90%  "
Set-TextValue "D9" "0.07181"
$ws.Range("E9").Value = "  +0.22%  "
Set-TextValue "D10" "24.67"
$ws.Range("E10").Value = "  -4.11%  "
Set-TextValue "D11" "0.08363"
$ws.Range("E11").Value = "  -2.97%  "
Set-TextValue "D12" "0.7498"
$ws.Range("E12").Value = "  -2.98%  "

# Rows 13/14 swap ranking: Polkadot <-> WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.891.31"
$ws.Range("E13").Value = "  -8.41%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "5.413"
$ws.Range("E14").Value = "  +0.56%  "

Set-TextValue "D15" "92.53"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "29.887.43"
$ws.Range("E16").Value = "  -1.06%  "
Set-TextValue "D17" "6.066"
$ws.Range("E17").Value = "  -1.86%  "
Set-TextValue "D18" "246.11"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E19").Value = "  -2.30%  "
Set-TextValue "D20" "0.000007821"
$ws.Range("E20").Value = "  +0.03%  "
Set-TextValue "D21" "0.9981"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "2.125.11"
$ws.Range("E22").Value = "  -9.64%  "
Set-TextValue "D23" "8.002"
$ws.Range("E23").Value = "  -0.57%  "
Set-TextValue "D24" "0.9979"
$ws.Range("E24").Value = "  -0.37%  "
Set-TextValue "D25" "0.1548"
$ws.Range("E25").Value = "  -5.82%  "
Set-TextValue "D26" "9.249"
$ws.Range("E26").Value = "  -1.56%  "
Set-TextValue "D27" "164.74"
$ws.Range("E27").Value = "  +1.20%  "
Set-TextValue "D28" "18.65"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("E29").Value = "  -0.52%  "
Set-TextValue "D30" "1.508"
$ws.Range("E30").Value = "  +5.01%  "
Set-TextValue "D31" "4.588"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("E32").Value = "  -0.67%  "
Set-TextValue "D33" "4.292"
$ws.Range("E33").Value = "  +4.05%  "
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("E35").Value = "  -0.46%  "
Set-TextValue "D36" "0.7543"
$ws.Range("E36").Value = "  +0.22%  "
Set-TextValue "D37" "0.9979"
$ws.Range("E37").Value = "  -0.84%  "
Set-TextValue "D38" "2.688"
$ws.Range("E38").Value = "  -0.37%  "
Set-TextValue "D39" "0.01959"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("E40").Value = "  -1.28%  "
Set-TextValue "D41" "0.4502"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "1.111.20"
Set-TextValue "D43" "6.050"
$ws.Range("E43").Value = "  -0.80%  "
Set-TextValue "D44" "72.28"
$ws.Range("E44").Value = "  -2.21%  "
Set-TextValue "D45" "0.8560"
$ws.Range("E45").Value = "  +0.61%  "

# Rows 46/47 swap ranking: Quant <-> PaxDollar
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D46" "1.000"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D47" "102.71"
$ws.Range("E47").Value = "  -1.16%  "

Set-TextValue "D48" "7.611"
$ws.Range("E48").Value = "  -0.08%  "

# Rows 49/50 swap ranking: SynthetixNetwork <-> RenderToken
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "1.841"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue "D50" "3.023"
$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("D51").Value = "2.021.64"
$ws.Range("E51").Value = "  -9.66%  "
